# Adds the "ECS + Atlas" worksheet with an hourly/daily/monthly AWS cost
# breakdown (ECR storage, ECS vCPU/GB, VPC IP/PrivateLink), and tidies up
# a few leftover cell-format quirks on the "App runner" sheet plus a new
# NAT-gateway reference block there.

$wb = $excel.ActiveWorkbook
$appRunner = $wb.Worksheets.Item("App runner")

# --- "App runner" sheet: small reference block for the NAT gateway cost ---
$appRunner.Range("G4").Value2 = "AWS Nat gateway"
$appRunner.Range("H4").Value2 = "per hour"
$appRunner.Range("I4").Value2 = 0.05
$appRunner.Columns.Item(7).ColumnWidth = 15.75

# --- "App runner" sheet: drop the stray fill/border variants that used to
#     sit on A13/A14/A15/A16/A17 (style got simplified back to defaults /
#     the plain bottom-border style already used elsewhere on the sheet) ---
$appRunner.Cells.Item(13,1).ClearFormats()
$appRunner.Cells.Item(15,1).ClearFormats()
$appRunner.Cells.Item(16,1).ClearFormats()

$appRunner.Cells.Item(14,1).Borders.Item(9).LineStyle = 1
$appRunner.Cells.Item(14,1).Borders.Item(9).ColorIndex = 64
$appRunner.Cells.Item(17,1).Borders.Item(9).LineStyle = 1
$appRunner.Cells.Item(17,1).Borders.Item(9).ColorIndex = 64

# the selection / active tab had moved off "App runner" onto the new sheet
[void]$appRunner.Range("G9").Select()

# --- add the new "ECS + Atlas" sheet right after "App runner" ---
$ecs = $wb.Worksheets.Add($null, $appRunner)
$ecs.Name = "ECS + Atlas"

# Header row
$ecs.Range("A1").Value2 = "Service"
$ecs.Range("B1").Value2 = "Element"
$ecs.Range("C1").Value2 = "per hour"
$ecs.Range("D1").Value2 = "units"
$ecs.Range("E1").Value2 = "total/hour"
$ecs.Range("F1").Value2 = "Hour cost"
$ecs.Range("G1").Value2 = "Day Cost"
$ecs.Range("G1").NumberFormat = "0.00"
$ecs.Range("H1").Value2 = "Month cost"

# ECR block
$ecs.Range("A2").Value2 = "ECR"
$ecs.Range("G2").NumberFormat = "0.00"

$ecs.Range("B3").Value2 = "Storage GB/hour"
$ecs.Range("C3").Formula = "=0.1/(24*30)"
$ecs.Range("D3").Value2 = 1
$ecs.Range("E3").Formula = "=D3*C3"
$ecs.Range("G3").NumberFormat = "0.00"

$ecs.Range("F4").Formula = "=E3"
$ecs.Range("G4").Formula = "=F4*24"
$ecs.Range("H4").Formula = "=G4*31"
$ecs.Range("F4:H4").NumberFormat = "0.00"

# ECS block
$ecs.Range("A5").Value2 = "ECS"
$ecs.Range("F5:H5").NumberFormat = "0.00"

$ecs.Range("B6").Value2 = "vCPU"
$ecs.Range("C6").Value2 = 0.04
$ecs.Range("D6").Value2 = 1
$ecs.Range("E6").Formula = "=D6*C6"
$ecs.Range("F6:H6").NumberFormat = "0.00"

$ecs.Range("B7").Value2 = "GB"
$ecs.Range("C7").Value2 = 0.0044
$ecs.Range("D7").Value2 = 2
$ecs.Range("E7").Formula = "=D7*C7"
$ecs.Range("F7:H7").NumberFormat = "0.00"

$ecs.Range("F8").Formula = "=SUM(E6:E7)"
$ecs.Range("G8").Formula = "=F8*24"
$ecs.Range("H8").Formula = "=G8*31"
$ecs.Range("F8:H8").NumberFormat = "0.00"

# VPC block
$ecs.Range("A9").Value2 = "VPC"
$ecs.Range("F9:H9").NumberFormat = "0.00"

$ecs.Range("B10").Value2 = "IP address"
$ecs.Range("C10").Value2 = 0.005
$ecs.Range("D10").Value2 = 1
$ecs.Range("E10").Formula = "=D10*C10"
$ecs.Range("F10:H10").NumberFormat = "0.00"

$ecs.Range("B11").Value2 = "Private link"
$ecs.Range("C11").Value2 = 0.011
$ecs.Range("D11").Value2 = 1
$ecs.Range("E11").Formula = "=D11*C11"
$ecs.Range("F11:H11").NumberFormat = "0.00"

$ecs.Range("F12").Formula = "=SUM(E10:E11)"
$ecs.Range("G12").Formula = "=F12*24"
$ecs.Range("H12").Formula = "=G12*31"
$ecs.Range("F12:H12").NumberFormat = "0.00"

$ecs.Range("F13:G13").NumberFormat = "0.00"
$ecs.Range("F14:G14").NumberFormat = "0.00"

# Totals
$ecs.Range("A15").Value2 = "Total"
$ecs.Range("F15").Formula = "=SUM(F1:F14)"
$ecs.Range("G15").Formula = "=SUM(G1:G14)"
$ecs.Range("H15").Formula = "=SUM(H1:H14)"
$ecs.Range("F15:H15").NumberFormat = "0.00"

$ecs.Range("F16").Value2 = "Hour"
$ecs.Range("G16").Value2 = "Day"
$ecs.Range("H16").Value2 = "Month"

# layout
$ecs.Columns.Item(1).ColumnWidth = 24.6
$ecs.Columns.Item(2).ColumnWidth = 15.6
$ecs.Columns.Item(5).ColumnWidth = 9.25
$ecs.Columns.Item(6).ColumnWidth = 8.6

[void]$ecs.Range("D6").Select()
[void]$ecs.Activate()
